# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The monthly "Periodo Mora" (period) / "Valor Mora" (amount) rows (16-27)
# get re-ordered: the list of periods 1701..1712 is reversed (so row 16 now
# shows period 1712, row 27 shows period 1701), and the "Valor Mora" amounts
# follow each period to its new row (only the 1712 row carries 12039, every
# other period keeps 32834).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periods in their new (reversed) top-to-bottom order, rows 16 through 27.
$periods = @("1712","1711","1710","1709","1708","1707","1706","1705","1704","1703","1702","1701")

# Valor Mora amount that travels together with each period.
$amounts = @{
    "1701" = 32834
    "1702" = 32834
    "1703" = 32834
    "1704" = 32834
    "1705" = 32834
    "1706" = 32834
    "1707" = 32834
    "1708" = 32834
    "1709" = 32834
    "1710" = 32834
    "1711" = 32834
    "1712" = 12039
}

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $period = $periods[$i]
    $ws.Range("E$row").Value = $period
    $ws.Range("F$row").Value = $amounts[$period]
}
